# Weekly data update: add the Wk28 entries to the "Weekly Expenditure" sheet.
#
# The four new rows (WorkWeek=Wk28) are exact duplicates of the four most
# recent existing rows (WorkWeek=Wk27, rows 46-49) except for the WorkWeek
# label itself - same Type/Material Number/Description/Quantity/SYSTEM/
# Date/Cost values. We duplicate them with Range.Copy so that cell typing
# (e.g. the Material Number in row 48/52 is a real number while the others
# are text-that-looks-numeric) and cell styles carry over exactly, then
# overwrite column A on the new rows with the new week label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRows = $ws.Range("A46:H49")
$targetRows = $ws.Range("A50:H53")
$sourceRows.Copy($targetRows) | Out-Null

$ws.Range("A50:A53").Value = "Wk28"

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("B48").Select() | Out-Null
